# Auto-generated Excel COM-interop edit script
# Applies updated market-price / profit values (H:N) for specific leve rows
# across multiple worksheets, per the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 351.7143
$ws.Range("I38").Value = 230.6
$ws.Range("J38").Value = 654.5
$ws.Range("K38").Value = 691.8
$ws.Range("L38").Value = 1963.5
$ws.Range("M38").Value = -319.8
$ws.Range("N38").Value = -2707.5

$ws.Range("H42").Value = 40.25
$ws.Range("I42").Value = 31
$ws.Range("J42").Value = 68
$ws.Range("K42").Value = 93
$ws.Range("L42").Value = 204
$ws.Range("M42").Value = 137
$ws.Range("N42").Value = -664

$ws.Range("H112").Value = 15626082
$ws.Range("J112").Value = 15626082
$ws.Range("L112").Value = 46878246
$ws.Range("N112").Value = -46880462

$ws.Range("H113").Value = 17603.5
$ws.Range("I113").Value = 22804.666
$ws.Range("K113").Value = 22804.666
$ws.Range("M113").Value = -19550.666

$ws.Range("H129").Value = 888.4605
$ws.Range("J129").Value = 936.19403
$ws.Range("L129").Value = 2808.58209
$ws.Range("N129").Value = -12808.58209

$ws.Range("H138").Value = 6415214.5
$ws.Range("I138").Value = 2017252.5
$ws.Range("J138").Value = 7465474
$ws.Range("K138").Value = 6051757.5
$ws.Range("L138").Value = 22396422
$ws.Range("M138").Value = -6046617.5
$ws.Range("N138").Value = -22406702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 350
$ws.Range("J17").Value = 350
$ws.Range("L17").Value = 350
$ws.Range("N17").Value = -696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2256.5557
$ws.Range("I94").Value = 2288.625
$ws.Range("K94").Value = 2288.625
$ws.Range("M94").Value = -1837.625

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45199.6
$ws.Range("J20").Value = 45199.6
$ws.Range("L20").Value = 45199.6
$ws.Range("N20").Value = -45671.6

$ws.Range("H30").Value = 45199.6
$ws.Range("J30").Value = 45199.6
$ws.Range("L30").Value = 45199.6
$ws.Range("N30").Value = -45381.6

$ws.Range("H31").Value = 1833.6818
$ws.Range("I31").Value = 1322.05
$ws.Range("K31").Value = 1322.05
$ws.Range("M31").Value = -1027.05

$ws.Range("H34").Value = 1833.6818
$ws.Range("I34").Value = 1322.05
$ws.Range("K34").Value = 1322.05
$ws.Range("M34").Value = -1120.05

$ws.Range("H128").Value = 45199.6
$ws.Range("J128").Value = 45199.6
$ws.Range("L128").Value = 45199.6
$ws.Range("N128").Value = -55159.6

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H135").Value = 19961.666
$ws.Range("J135").Value = 19961.666
$ws.Range("L135").Value = 19961.666
$ws.Range("N135").Value = -30101.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1495.7142
$ws.Range("I5").Value = 715.3077
$ws.Range("J5").Value = 2763.875
$ws.Range("K5").Value = 2145.9231
$ws.Range("L5").Value = 8291.625
$ws.Range("M5").Value = -2033.9231
$ws.Range("N5").Value = -8515.625

$ws.Range("H22").Value = 1336
$ws.Range("I22").Value = 225
$ws.Range("J22").Value = 1780.4
$ws.Range("K22").Value = 675
$ws.Range("L22").Value = 5341.200000000001
$ws.Range("M22").Value = -506
$ws.Range("N22").Value = -5679.200000000001

$ws.Range("H27").Value = 1336
$ws.Range("I27").Value = 225
$ws.Range("J27").Value = 1780.4
$ws.Range("K27").Value = 675
$ws.Range("L27").Value = 5341.200000000001
$ws.Range("M27").Value = -573
$ws.Range("N27").Value = -5545.200000000001

$ws.Range("H34").Value = 2064.182
$ws.Range("I34").Value = 201
$ws.Range("J34").Value = 2762.875
$ws.Range("K34").Value = 603
$ws.Range("L34").Value = 8288.625
$ws.Range("M34").Value = -519
$ws.Range("N34").Value = -8456.625

$ws.Range("H58").Value = 1901.25
$ws.Range("I58").Value = 1901.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5703.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5575.75
$ws.Range("N58").ClearContents()

$ws.Range("H112").Value = 166668380
$ws.Range("I112").Value = 1130.6
$ws.Range("J112").Value = 285716400
$ws.Range("K112").Value = 3391.8
$ws.Range("L112").Value = 857149200
$ws.Range("M112").Value = -2283.8
$ws.Range("N112").Value = -857151416

$ws.Range("H113").Value = 970.2083
$ws.Range("J113").Value = 986.3043
$ws.Range("L113").Value = 2958.9129
$ws.Range("N113").Value = -7298.9129

$ws.Range("H132").Value = 1163.7894
$ws.Range("I132").Value = 839
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 7551
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -5021
$ws.Range("N132").Value = -17660

$ws.Range("H135").Value = 1495.7142
$ws.Range("I135").Value = 715.3077
$ws.Range("J135").Value = 2763.875
$ws.Range("K135").Value = 6437.7693
$ws.Range("L135").Value = 24874.875
$ws.Range("M135").Value = -3902.7693
$ws.Range("N135").Value = -29944.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1252.5
$ws.Range("I97").Value = 1174.2858
$ws.Range("K97").Value = 1174.2858
$ws.Range("M97").Value = -678.2858000000001

$ws.Range("H103").Value = 18666.666
$ws.Range("J103").Value = 18666.666
$ws.Range("L103").Value = 18666.666
$ws.Range("N103").Value = -21010.666

$ws.Range("H107").Value = 416.7
$ws.Range("J107").Value = 549.5
$ws.Range("L107").Value = 549.5
$ws.Range("N107").Value = -4389.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1116
$ws.Range("I46").Value = 1049.6666
$ws.Range("J46").Value = 1160.2222
$ws.Range("K46").Value = 1049.6666
$ws.Range("L46").Value = 1160.2222
$ws.Range("M46").Value = -861.6666
$ws.Range("N46").Value = -1536.2222

$ws.Range("H132").Value = 3478.32
$ws.Range("I132").Value = 2294.6667
$ws.Range("K132").Value = 6884.000100000001
$ws.Range("M132").Value = -4354.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 530.93335
$ws.Range("I107").Value = 496.5
$ws.Range("K107").Value = 1489.5
$ws.Range("M107").Value = 430.5

$ws.Range("H132").Value = 26322234
$ws.Range("I132").Value = 50008410
$ws.Range("K132").Value = 150025230
$ws.Range("M132").Value = -150022700

$ws.Range("H136").Value = 17598550
$ws.Range("I136").Value = 25718942
$ws.Range("J136").Value = 4369.5
$ws.Range("K136").Value = 77156826
$ws.Range("L136").Value = 13108.5
$ws.Range("M136").Value = -77154276
